# Force Mastercard Orange (#FF5F00) on all heading text.
#
# 1) Every Heading1/Heading2/Heading3 paragraph in the body gets an explicit
#    run-level <w:color w:val="FF5F00"/> on its (single) run, so the
#    orange shows even though the style still carries a theme color.
# 2) The paragraph styles Heading5, Heading6, Title and Subtitle are
#    redefined to use the flat RGB color FF5F00 instead of a theme color,
#    so any theme swap can no longer override them.

$d = $word.ActiveDocument

# wdColorAutomatic-style literal: Word VBA RGB() packs as 0x00BBGGRR.
# FF5F00 -> R=0xFF G=0x5F B=0x00 -> 0x005FFF = 24575
$mastercardOrange = 24575

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    if ($styleName -match "^Heading [1-3]$") {
        $r = $p.Range
        # Exclude the trailing paragraph mark so only the run(s) holding the
        # heading text get the explicit color (not the pilcrow / pPr rPr).
        [void]$r.MoveEnd(1, -1)
        $r.Font.Color = $mastercardOrange
    }
}

foreach ($styleName in @("Heading5", "Heading6", "Title", "Subtitle")) {
    $s = $d.Styles($styleName)
    $s.Font.Color = $mastercardOrange
}
